$d = $word.ActiveDocument

# Locate "This is a Microsoft word document." and append " (Changed main)"
# right after it (before the paragraph mark), matching the target diff
# which adds three extra runs: " (", "Changed main", ")".
$rng = $d.Content
$found = $rng.Find.Execute("This is a Microsoft word document.", $true, $false, `
                            $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Collapse(0)          # wdCollapseEnd - move to just after the found text
    $rng.InsertAfter(" (")
    $rng.Collapse(0)
    $rng.InsertAfter("Changed main")
    $rng.Collapse(0)
    $rng.InsertAfter(")")
}
